# Auto update Excel log
# Appends newly-logged sensor events to the SeniorConnect master log.
#
#  - Proximity        : 3 new "Bedroom Door" enter/exit events (rows 2-4)
#  - mmWave(BR)        : 1 new "Empty" reading, numeric value 0 (row 14)
#  - mmWave(HR)        : 1 new "Empty" reading, numeric value 0 (row 14)
#  - mmWave(InBed)     : 1 new "Out of Bed" / "Empty" reading (row 14)

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $addr, $value) {
    # Force text storage so date-like strings (e.g. "2026-02-01") are not
    # silently reinterpreted as Excel date serials.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
}

# ---------------------------------------------------------------------
# Proximity - Bedroom Door sensor: three new enter/exit rows
# ---------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")

Set-TextCell $wsProximity "A2" "2026-02-01"
$wsProximity.Range("B2").Value = "15:04:14"
$wsProximity.Range("C2").Value = "15:00"
$wsProximity.Range("D2").Value = "Bedroom Door"
$wsProximity.Range("E2").Value = "ENTER"
$wsProximity.Range("F2").Value = "User ENTERED Bedroom"

Set-TextCell $wsProximity "A3" "2026-02-01"
$wsProximity.Range("B3").Value = "15:04:21"
$wsProximity.Range("C3").Value = "15:00"
$wsProximity.Range("D3").Value = "Bedroom Door"
$wsProximity.Range("E3").Value = "EXIT"
$wsProximity.Range("F3").Value = "User EXITED Bedroom"

Set-TextCell $wsProximity "A4" "2026-02-01"
$wsProximity.Range("B4").Value = "15:04:34"
$wsProximity.Range("C4").Value = "15:00"
$wsProximity.Range("D4").Value = "Bedroom Door"
$wsProximity.Range("E4").Value = "ENTER"
$wsProximity.Range("F4").Value = "User ENTERED Bedroom"

# ---------------------------------------------------------------------
# mmWave(BR) - new "Empty" reading (numeric value 0)
# ---------------------------------------------------------------------
$wsBR = $wb.Worksheets.Item("mmWave(BR)")

Set-TextCell $wsBR "A14" "2026-02-01"
$wsBR.Range("B14").Value = "15:04:34"
$wsBR.Range("C14").Value = "15:00"
$wsBR.Range("D14").Value = "Bedroom"
$wsBR.Range("E14").Value = 0
$wsBR.Range("F14").Value = "Empty"

# ---------------------------------------------------------------------
# mmWave(HR) - new "Empty" reading (numeric value 0)
# ---------------------------------------------------------------------
$wsHR = $wb.Worksheets.Item("mmWave(HR)")

Set-TextCell $wsHR "A14" "2026-02-01"
$wsHR.Range("B14").Value = "15:04:34"
$wsHR.Range("C14").Value = "15:00"
$wsHR.Range("D14").Value = "Bedroom"
$wsHR.Range("E14").Value = 0
$wsHR.Range("F14").Value = "Empty"

# ---------------------------------------------------------------------
# mmWave(InBed) - new "Out of Bed" / "Empty" reading
# ---------------------------------------------------------------------
$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")

Set-TextCell $wsInBed "A14" "2026-02-01"
$wsInBed.Range("B14").Value = "15:04:33"
$wsInBed.Range("C14").Value = "15:00"
$wsInBed.Range("D14").Value = "Bedroom"
$wsInBed.Range("E14").Value = "Out of Bed"
$wsInBed.Range("F14").Value = "Empty"
